# Add new "2022-Q4" sheet as a copy of "2022-Q3" (so it inherits identical
# layout/styling), inserted right before "2022-Q3" (i.e. right after "总计").
$wb = $excel.ActiveWorkbook
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Fill in the new quarter's fund data on the "2022-Q4" sheet.
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "539003"
$q4.Range("C2").Value = "建信富时100指数（QDII）人民币A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.56"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "85.06"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "5.49"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0307"
$q4.Range("H2").Value = 5

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "008707"
$q4.Range("C3").Value = "建信富时100指数（QDII）美元现汇 A"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.56"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "85.06"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "5.49"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0307"
$q4.Range("H3").Value = 5

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "008706"
$q4.Range("C4").Value = "建信富时100指数（QDII）人民币 C"
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "0.26"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "85.06"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "5.49"
$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "0.0143"
$q4.Range("H4").Value = 5

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "008708"
$q4.Range("C5").Value = "建信富时100指数（QDII）美元现汇 C"
$q4.Range("D5").NumberFormat = "@"
$q4.Range("D5").Value = "0.26"
$q4.Range("E5").NumberFormat = "@"
$q4.Range("E5").Value = "85.06"
$q4.Range("F5").NumberFormat = "@"
$q4.Range("F5").Value = "5.49"
$q4.Range("G5").NumberFormat = "@"
$q4.Range("G5").Value = "0.0143"
$q4.Range("H5").Value = 5

# Update the "总计" summary sheet: insert a new row for 2022-Q4 right after
# the header, shifting the existing quarter rows down by one.
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 4
$zj.Range("D2").Value = 0.09

For ($i = 3; $i -le 10; $i++) {
    $zj.Cells.Item($i, 1).Value = $i - 2
}

# Restore formatting on the newly inserted row to match the other data rows.
$zj.Range("A3:D3").Copy()
$zj.Range("A2:D2").PasteSpecial(-4122)

# The last tab ("2020-Q4") should remain the active/selected sheet.
$last = $wb.Worksheets.Item("2020-Q4")
$last.Activate()
